$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column A (rows 10-13) and B17
$ws.Range("A10").Value = "a"
$ws.Range("A11").Value = "awd"
$ws.Range("A12").Value = "dwa"
$ws.Range("A13").Value = "d"

# Add new values in column A (rows 14-16)
$ws.Range("A14").Value = "awd"
$ws.Range("A15").Value = "dad"
$ws.Range("A16").Value = "awd"

# Update B17 value
$ws.Range("B17").Value = "d"

# Update the selected cell to B17
$null = $ws.Range("B17").Select()
